$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Cypher query for the "ParticipantsTab" row (B2) - the query now
# walks participant -> study / sample with OPTIONAL MATCHes (so rows with
# no diagnosis/sample/genomic_info are still returned), re-collects the
# samples after filtering, and sorts them.
$newParticipantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['JSON']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$ws.Range("B2").Value2 = $newParticipantQuery

# Update the view: scrolled position moved so row 2 is the top visible row
# (was row 4), and the selected cell moved from C5 to B4.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("B4").Select()
